# edit.ps1
# Rewrites the "Embracing Complexity" essay into "The Intricate World of
# Chemistry", changes the author's name/email, and appends a trailing
# empty paragraph, per the target diff.

$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $r = $d.Content
    $r.Find.ClearFormatting()
    $ok = $r.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
    if (-not $ok) {
        throw "Find/Replace failed for: $old"
    }
}

function Split-At($pos) {
    # Force a run boundary right before character at $pos by toggling a
    # character formatting property on/off for the single character at
    # [pos, pos+1). This leaves formatting unchanged but breaks the run
    # into two (or three) separate <w:r> elements when saved.
    $rr = $d.Range($pos, $pos + 1)
    $rr.Font.Bold = 1
    $rr.Font.Bold = 0
}

function Split-After($text) {
    # Find the (already-replaced) $text and make sure whatever follows it
    # starts in its own run.
    $r = $d.Content
    $r.Find.ClearFormatting()
    $ok = $r.Find.Execute($text, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
    if (-not $ok) {
        throw "Split-After: could not find: $text"
    }
    Split-At($r.End)
}

# =================================================================
# Title
# =================================================================
Replace-Text "Embracing Complexity: The Path to Profound Understanding" "The Intricate World of Chemistry: Unveiling the Hidden Symphony of Elements and Molecules"

# =================================================================
# Author name: "Dr. Emily Foster" (3 runs) -> "Sofia Rodriguez" (1 run)
# =================================================================
Replace-Text "Dr. Emily Foster" "Sofia Rodriguez"

# =================================================================
# Email: "emily_foster@knowledgehub" -> "sofia" + "." + "rodriguez",
# keeping the existing "." and "edu" runs intact and separate.
# =================================================================
Replace-Text "emily_foster@knowledgehub" "sofia.rodriguez"
Split-After "sofia"
Split-After "sofia.rodriguez"

# =================================================================
# Body paragraph 1 (intro, 4 sentences before the first <br/><br/>)
# =================================================================
Replace-Text "In an era of rapid advancements and specialized knowledge, we often find ourselves thrust into the depths of complexity" "Chemistry, a captivating branch of science, unravels the profound mysteries of matter and its transformations"
Replace-Text " From the intricacies of quantum mechanics to the interconnectedness of global ecosystems, understanding the world around us requires a willingness to navigate the inherent complexities that exist within every field of study" " From the grand symphony of chemical reactions to the intricate dance of atoms, chemistry unveils the principles governing our material world"
Replace-Text " The path to profound understanding lies not in simplifying or reducing the complexity of our inquiries, but in embracing and exploring it" " At the core of chemistry lies the manipulation of elements, those fundamental building blocks that combine to form all substances"
Replace-Text " This essay delves into the reasons why complexity should be cherished in the pursuit of knowledge and the transformative insights it can offer across diverse disciplines" " By delving into the depths of chemical structure and reactivity, we gain insight into the composition and properties of matter, laying the groundwork for comprehending the intricate tapestry of life and the universe"

Split-After "Chemistry, a captivating branch of science, unravels the profound mysteries of matter and its transformations"
Split-After " From the grand symphony of chemical reactions to the intricate dance of atoms, chemistry unveils the principles governing our material world"
Split-After " At the core of chemistry lies the manipulation of elements, those fundamental building blocks that combine to form all substances"
Split-After " By delving into the depths of chemical structure and reactivity, we gain insight into the composition and properties of matter, laying the groundwork for comprehending the intricate tapestry of life and the universe"

# =================================================================
# Body paragraph 1, second chunk (after first <br/><br/>)
# =================================================================
Replace-Text "Complexity holds the key to uncovering hidden relationships, patterns, and correlations" "Unveiling the secrets of chemistry leads us to explore the profound interconnections between elements and their dynamic interactions"
Replace-Text " By delving into its depths, scientists can uncover the interdependencies between seemingly disparate phenomena" " We uncover the patterns governing chemical bonding, the forces that hold atoms together, creating substances with unique properties and functions"
Replace-Text " In economics, the study of complex systems reveals how individual decisions and actions within a market can lead to unpredictable market behaviors" " This knowledge forms the foundation for understanding materials science, enabling us to develop novel materials with tailored properties and applications"

# The ecology / medicine / politics sentences (3 sentences + 2 "." runs)
# collapse into a single new sentence.
Replace-Text " In ecology, understanding the intricate interdependencies among species unveils the delicate balance that sustains ecosystems. In medicine, mapping the complex interplay of genes within a biological system can lead to breakthroughs in disease diagnosis and treatment. In politics, grasping the multifaceted dynamics of international relations can empower policymakers with the knowledge to avert conflicts and promote global peace" " Furthermore, by unraveling the intricate dance of chemical reactions, we harness their energy to power our world, providing us with electricity, heat, and fuels"

Split-After "Unveiling the secrets of chemistry leads us to explore the profound interconnections between elements and their dynamic interactions"
Split-After " We uncover the patterns governing chemical bonding, the forces that hold atoms together, creating substances with unique properties and functions"
Split-After " This knowledge forms the foundation for understanding materials science, enabling us to develop novel materials with tailored properties and applications"
Split-After " Furthermore, by unraveling the intricate dance of chemical reactions, we harness their energy to power our world, providing us with electricity, heat, and fuels"

# =================================================================
# Body paragraph 1, third chunk (after second <br/><br/>)
# =================================================================
Replace-Text "Furthermore, complexity demands a multidisciplinary approach, fostering collaboration between scholars from various fields" "Chemistry plays a pivotal role in unlocking the enigma of biological processes, the symphony of life"
Replace-Text " By pooling their knowledge and perspectives, researchers can tackle complex problems that lie beyond the grasp of any single discipline" " By examining the intricate interactions of biomolecules, we gain insights into the mechanisms underlying metabolism, genetics, and disease"
Replace-Text " Physicists engage with biologists to explore the fundamental building blocks of life" " This understanding has led to breakthroughs in medicine, the development of life-saving drugs, and advancements in gene therapy"

# The computer-scientists / neurologists / "through these collaborations" /
# lastRenderedPageBreak / "diverse perspectives" sentences all collapse
# into one new sentence; the lastRenderedPageBreak marker disappears here
# (it reappears, fresh, in the Summary paragraph below).
Replace-Text " Computer scientists join forces with archaeologists to uncover ancient secrets hidden in digital artifacts. Neurologists collaborate with musicians to study the neural basis of creativity. Through these collaborations, diverse perspectives intersect, leading to innovative solutions and a deeper understanding of the interconnectedness of all things" " Moreover, chemistry empowers us to unlock the secrets of nutrition, enabling us to comprehend the intricate relationship between diet and health"

Split-After "Chemistry plays a pivotal role in unlocking the enigma of biological processes, the symphony of life"
Split-After " By examining the intricate interactions of biomolecules, we gain insights into the mechanisms underlying metabolism, genetics, and disease"
Split-After " This understanding has led to breakthroughs in medicine, the development of life-saving drugs, and advancements in gene therapy"
Split-After " Moreover, chemistry empowers us to unlock the secrets of nutrition, enabling us to comprehend the intricate relationship between diet and health"

Write-Host "stage3 (para4 second half) done"
